# Weekly update: prepend a new Cilantro price observation for
# Terminal Hortofrutícola Agro Chillán (row 8), pushing the existing
# rows 8:27 down to 9:28. The sheet's used range grows from A1:R27 to
# A1:R28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8 - this shifts the former rows 8:27
# down to 9:28 (data, formatting and dimension all update automatically).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with this week's record.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44623
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 550
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = 575
$ws.Range("N8").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O8").Value = "Provincia de Diguillín"
$ws.Range("P8").Value = 575
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
